$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H103").Value = 141.77777
$ws.Range("J103").Value = 155.75
$ws.Range("L103").Value = 467.25
$ws.Range("N103").Value = -1639.25
$ws.Range("H106").Value = 6946900.5
$ws.Range("I106").Value = 12822188
$ws.Range("J106").Value = 3377.9092
$ws.Range("K106").Value = 12822188
$ws.Range("L106").Value = 3377.9092
$ws.Range("M106").Value = -12821557
$ws.Range("N106").Value = -4639.9092
$ws.Range("H116").Value = 5792.4707
$ws.Range("I116").Value = 3994
$ws.Range("J116").Value = 6541.8335
$ws.Range("K116").Value = 3994
$ws.Range("L116").Value = 6541.8335
$ws.Range("M116").Value = -552
$ws.Range("N116").Value = -13425.8335
$ws.Range("H137").Value = 1540.6364
$ws.Range("I137").Value = 1280.0555
$ws.Range("K137").Value = 3840.1665
$ws.Range("M137").Value = -1290.1665
$ws.Range("H138").Value = 33336660
$ws.Range("I138").Value = 142860400
$ws.Range("J138").Value = 3348.5217
$ws.Range("K138").Value = 428581200
$ws.Range("L138").Value = 10045.5651
$ws.Range("M138").Value = -428576060
$ws.Range("N138").Value = -20325.5651
$ws.Range("H141").Value = 3228.3845
$ws.Range("I141").Value = 2656.9
$ws.Range("K141").Value = 7970.700000000001
$ws.Range("M141").Value = -2790.700000000001

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H32").Value = 6125.582
$ws.Range("I32").Value = 5292.0195
$ws.Range("K32").Value = 5292.0195
$ws.Range("M32").Value = -5005.0195
$ws.Range("H43").Value = 0
$ws.Range("J43").Value = 0
$ws.Range("N43").Value = 0
$ws.Range("L43").ClearContents()
$ws.Range("H45").Value = 2166.238
$ws.Range("I45").Value = 1631.7273
$ws.Range("J45").Value = 2754.2
$ws.Range("K45").Value = 1631.7273
$ws.Range("L45").Value = 2754.2
$ws.Range("M45").Value = -1254.7273
$ws.Range("N45").Value = -3508.2
$ws.Range("H112").Value = 31433.5
$ws.Range("J112").Value = 31433.5
$ws.Range("L112").Value = 31433.5
$ws.Range("N112").Value = -34387.5
$ws.Range("H114").Value = 31487.125
$ws.Range("J114").Value = 31487.125
$ws.Range("L114").Value = 31487.125
$ws.Range("N114").Value = -40165.125
$ws.Range("H122").Value = 2052.8518
$ws.Range("I122").Value = 1323
$ws.Range("K122").Value = 3969
$ws.Range("M122").Value = -1519
$ws.Range("H125").Value = 31626.666
$ws.Range("J125").Value = 31626.666
$ws.Range("L125").Value = 31626.666
$ws.Range("N125").Value = -41466.666

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H20").Value = 2576
$ws.Range("I20").Value = 2366.3333
$ws.Range("J20").Value = 2995.3333
$ws.Range("K20").Value = 2366.3333
$ws.Range("L20").Value = 2995.3333
$ws.Range("M20").Value = -2119.3333
$ws.Range("N20").Value = -3489.3333
$ws.Range("H55").Value = 0
$ws.Range("J55").Value = 0
$ws.Range("N55").Value = 0
$ws.Range("L55").ClearContents()
$ws.Range("H105").Value = 2084576
$ws.Range("I105").Value = 1065.1666
$ws.Range("J105").Value = 8335108.5
$ws.Range("K105").Value = 1065.1666
$ws.Range("L105").Value = 8335108.5
$ws.Range("M105").Value = 681.8334
$ws.Range("N105").Value = -8338602.5
$ws.Range("H134").Value = 2967.6047
$ws.Range("I134").Value = 3376.4722
$ws.Range("J134").Value = 864.8570999999999
$ws.Range("K134").Value = 10129.4166
$ws.Range("L134").Value = 2594.5713
$ws.Range("M134").Value = -7594.4166
$ws.Range("N134").Value = -7664.5713
$ws.Range("H135").Value = 29980
$ws.Range("J135").Value = 29980
$ws.Range("L135").Value = 29980
$ws.Range("N135").Value = -40120

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H31").Value = 3361.4167
$ws.Range("I31").Value = 2313.2632
$ws.Range("J31").Value = 4532.8823
$ws.Range("K31").Value = 2313.2632
$ws.Range("L31").Value = 4532.8823
$ws.Range("M31").Value = -2018.2632
$ws.Range("N31").Value = -5122.8823
$ws.Range("H34").Value = 3361.4167
$ws.Range("I34").Value = 2313.2632
$ws.Range("J34").Value = 4532.8823
$ws.Range("K34").Value = 2313.2632
$ws.Range("L34").Value = 4532.8823
$ws.Range("M34").Value = -2111.2632
$ws.Range("N34").Value = -4936.8823
$ws.Range("H99").Value = 17860118
$ws.Range("I99").Value = 2594.2856
$ws.Range("J99").Value = 35717640
$ws.Range("K99").Value = 2594.2856
$ws.Range("L99").Value = 35717640
$ws.Range("M99").Value = -1096.2856
$ws.Range("N99").Value = -35720636
$ws.Range("H122").Value = 1456.8
$ws.Range("I122").Value = 1485.7778
$ws.Range("J122").Value = 1413.3334
$ws.Range("K122").Value = 4457.3334
$ws.Range("L122").Value = 4240.0002
$ws.Range("M122").Value = -2007.3334
$ws.Range("N122").Value = -9140.0002
$ws.Range("H126").Value = 17860118
$ws.Range("I126").Value = 2594.2856
$ws.Range("J126").Value = 35717640
$ws.Range("K126").Value = 7782.8568
$ws.Range("L126").Value = 107152920
$ws.Range("M126").Value = -5312.8568
$ws.Range("N126").Value = -107157860

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H131").Value = 738.6882000000001
$ws.Range("I131").Value = 665
$ws.Range("J131").Value = 740.3077
$ws.Range("K131").Value = 1995
$ws.Range("L131").Value = 2220.9231
$ws.Range("M131").Value = 3045
$ws.Range("N131").Value = -12300.9231
$ws.Range("H138").Value = 1909.9231
$ws.Range("I138").Value = 1699.909
$ws.Range("J138").Value = 3065
$ws.Range("K138").Value = 5099.727000000001
$ws.Range("L138").Value = 9195
$ws.Range("M138").Value = 40.27299999999923
$ws.Range("N138").Value = -19475

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H70").Value = 12522440
$ws.Range("I70").Value = 3800
$ws.Range("J70").Value = 15652100
$ws.Range("K70").Value = 3800
$ws.Range("L70").Value = 15652100
$ws.Range("M70").Value = -3530
$ws.Range("N70").Value = -15652640
$ws.Range("H73").Value = 12522440
$ws.Range("I73").Value = 3800
$ws.Range("J73").Value = 15652100
$ws.Range("K73").Value = 3800
$ws.Range("L73").Value = 15652100
$ws.Range("M73").Value = -2864
$ws.Range("N73").Value = -15653972
$ws.Range("H97").Value = 1940.4546
$ws.Range("I97").Value = 1993.75
$ws.Range("J97").Value = 1876.5
$ws.Range("K97").Value = 1993.75
$ws.Range("L97").Value = 1876.5
$ws.Range("M97").Value = -1497.75
$ws.Range("N97").Value = -2868.5
$ws.Range("H126").Value = 5705
$ws.Range("I126").Value = 4684.615
$ws.Range("K126").Value = 14053.845
$ws.Range("M126").Value = -11583.845
$ws.Range("H132").Value = 38961.355
$ws.Range("I132").Value = 3881
$ws.Range("J132").Value = 58450.445
$ws.Range("K132").Value = 11643
$ws.Range("L132").Value = 175351.335
$ws.Range("M132").Value = -9113
$ws.Range("N132").Value = -180411.335

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H46").Value = 964.9545000000001
$ws.Range("I46").Value = 954.7646999999999
$ws.Range("J46").Value = 999.6
$ws.Range("K46").Value = 954.7646999999999
$ws.Range("L46").Value = 999.6
$ws.Range("M46").Value = -766.7646999999999
$ws.Range("N46").Value = -1375.6
$ws.Range("H100").Value = 1980.8823
$ws.Range("I100").Value = 1357.1428
$ws.Range("J100").Value = 2417.5
$ws.Range("K100").Value = 1357.1428
$ws.Range("L100").Value = 2417.5
$ws.Range("M100").Value = -816.1428000000001
$ws.Range("N100").Value = -3499.5
$ws.Range("H104").Value = 17815.5
$ws.Range("J104").Value = 17815.5
$ws.Range("L104").Value = 17815.5
$ws.Range("N104").Value = -24803.5
$ws.Range("H122").Value = 936763.0600000001
$ws.Range("I122").Value = 2453780.5
$ws.Range("J122").Value = 3213.8462
$ws.Range("K122").Value = 7361341.5
$ws.Range("L122").Value = 9641.5386
$ws.Range("M122").Value = -7358891.5
$ws.Range("N122").Value = -14541.5386
$ws.Range("H127").Value = 16713.312
$ws.Range("J127").Value = 16713.312
$ws.Range("L127").Value = 16713.312
$ws.Range("N127").Value = -26633.312
$ws.Range("H132").Value = 4705.857
$ws.Range("I132").Value = 4188.6
$ws.Range("K132").Value = 12565.8
$ws.Range("M132").Value = -10035.8

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H6").Value = 1740
$ws.Range("I6").Value = 1000
$ws.Range("J6").Value = 1925
$ws.Range("K6").Value = 1000
$ws.Range("L6").Value = 1925
$ws.Range("M6").Value = -885
$ws.Range("N6").Value = -2155
$ws.Range("H132").Value = 1578.9
$ws.Range("I132").Value = 1092.9412
$ws.Range("J132").Value = 4332.6665
$ws.Range("K132").Value = 3278.8236
$ws.Range("L132").Value = 12997.9995
$ws.Range("M132").Value = -748.8235999999997
$ws.Range("N132").Value = -18057.9995

